$d = $word.ActiveDocument

# --- Change 1: split "Types of wasteManagement and Legislation in Brazil" into two lines ---
$d.Content.Find.Execute(
    "Types of wasteManagement and Legislation in Brazil",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Types of waste" + [char]11 + "Management and Legislation in Brazil",
    2
) | Out-Null

# --- Change 2: break the Bibliografia paragraph into multiple lines via manual line breaks ---
$bibRange = $d.Content
$found = $bibRange.Find.Execute("Bibliografia básica:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bibRange.Expand(4) | Out-Null
    $bibRange.MoveEnd(1, -1) | Out-Null
    $bibRange.Text = 'Bibliografia básica:' + [char]11 + [char]11 + 'BARROS, R.M., Tratado sobre resíduos sólidos: gestão, uso e sustentabilidade, Editora Interciência, 2013.' + [char]11 + 'DA SILVA-FILHO, C.R.V., SOLER, F.D., Gestão de resíduos sólidos: o que diz a lei, 2° ed., Editora Trevisan, 2013.' + [char]11 + 'RIBEIRO, D.V., MORELLI, M.R., Resíduos sólidos: problemas ou oportunidades?, Editora Interciência, 2009.' + [char]11 + [char]11 + 'Bibliografia complementar:' + [char]11 + [char]11 + 'BRAGA B. (Org.), Introdução à engenharia ambiental: o desafio do desenvolvimento sustentável, 2° ed., Ed. Pearson Prentice Hall, 2005.' + [char]11 + 'CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão, Ed. Campus, 2013.' + [char]11 + 'CARVALHO, T.C.M.B., XAVIER, L.H. (Org.), Gestão de resíduos eletroeletrônicos: uma abordagem prática para a sustentabilidade, Edidora Elsivier Ltda, 2014.' + [char]11 + 'JACOBI, P. (Org.), Gestão compartilhada dos resíduos sólidos no Brasil: inovação com inclusão social, Annablume, 2006' + [char]11 + 'PEREIRA-NETO, J.T., Gerenciamento do lixo urbano: aspectos técnicos e operacionais, Editora UFV, 2013.' + [char]11 + 'SÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos, Ed. Oficina de textos, 3° reimpressão, 2011' + [char]11 + [char]11 + 'ASSOCIAÇÃO BRASILEIRA DE NORMAS TÉCNICAS (ABNT) – NORMAS ABNT RESÍDUOS SÓLIDOS: COLETÂNEA DE NORMAS: NBR 10.004, NBR 10.005, NBR 10.006, NBR 10.007' + [char]11 + 'ASSOCIAÇÃO BRSILEIRA DE NORMAS TÉCNICAS (ABNT) – NORMAS TÉCNICAS APRESENTAÇÃO DE PROJETOS DE ATERROS CONTROLADOS DE RESÍDUOS SÓLIDOS URBANOS: NBR 8849'
}
